$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(10).Cut()
$ws.Columns.Item(8).Insert()
$ws.Range("H1").Value = "Website"
